$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# "Date" value used in several places changes from 23 -> 24
$ws.Range("B12").Value = "24"
$ws.Range("B24").Value = "24"
$ws.Range("B35").Value = "24"
$ws.Range("B43").Value = "24"

# Formatted date value changes from 23-01-2022 -> 24-01-2022
$ws.Range("B15").Value = "24-01-2022"
# B29's existing format uses a quote-prefix (leading apostrophe) entry;
# re-supply it the same way so its cell formatting is preserved
$ws.Range("B29").Value = "'24-01-2022"
$ws.Range("B38").Value = "24-01-2022"
$ws.Range("B46").Value = "24-01-2022"

# Email value changes from nv31@gmail.com -> nv232@gmail.com
$ws.Range("B5").Value = "'nv232@gmail.com"
